$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of cell address -> new "Recorded By" value (System moved to front of the list)
$updates = @{
    'G2' = 'System, system, backup@backdoor.com'
    'G3' = 'System, dnasr281@gmail.com'
    'G5' = 'System, backup@backdoor.com'
    'G6' = 'System, dnasr281@gmail.com'
    'G7' = 'System, admin@admin.com'
    'G8' = 'System, backup@backdoor.com'
    'G10' = 'System, dnasr281@gmail.com'
    'G11' = 'System, dnasr281@gmail.com'
    'G12' = 'System, dnasr281@gmail.com'
    'G13' = 'System, dnasr281@gmail.com'
    'G14' = 'System, dnasr281@gmail.com'
    'G15' = 'System, dnasr281@gmail.com'
    'G17' = 'System, dnasr281@gmail.com'
    'G18' = 'System, dnasr281@gmail.com'
    'G19' = 'System, dnasr281@gmail.com'
    'G20' = 'System, dnasr281@gmail.com'
    'G21' = 'System, dnasr281@gmail.com'
    'G22' = 'System, dnasr281@gmail.com'
    'G24' = 'System, dnasr281@gmail.com'
    'G26' = 'System, dnasr281@gmail.com'
    'G28' = 'System, system, backup@backdoor.com'
    'G29' = 'System, dnasr281@gmail.com'
    'G31' = 'System, backup@backdoor.com'
    'G32' = 'System, dnasr281@gmail.com'
    'G33' = 'System, admin@admin.com'
    'G34' = 'System, backup@backdoor.com'
    'G36' = 'System, dnasr281@gmail.com'
    'G37' = 'System, dnasr281@gmail.com'
    'G38' = 'System, dnasr281@gmail.com'
    'G39' = 'System, dnasr281@gmail.com'
    'G40' = 'System, dnasr281@gmail.com'
    'G41' = 'System, dnasr281@gmail.com'
    'G43' = 'System, dnasr281@gmail.com'
    'G44' = 'System, dnasr281@gmail.com'
    'G45' = 'System, dnasr281@gmail.com'
    'G46' = 'System, dnasr281@gmail.com'
    'G47' = 'System, dnasr281@gmail.com'
    'G48' = 'System, dnasr281@gmail.com'
    'G50' = 'System, dnasr281@gmail.com'
    'G52' = 'System, dnasr281@gmail.com'
    'G54' = 'System, system, backup@backdoor.com'
    'G55' = 'System, dnasr281@gmail.com'
    'G57' = 'System, backup@backdoor.com'
    'G58' = 'System, dnasr281@gmail.com'
    'G59' = 'System, admin@admin.com'
    'G60' = 'System, backup@backdoor.com'
    'G62' = 'System, dnasr281@gmail.com'
    'G63' = 'System, dnasr281@gmail.com'
    'G64' = 'System, dnasr281@gmail.com'
    'G65' = 'System, dnasr281@gmail.com'
    'G66' = 'System, dnasr281@gmail.com'
    'G67' = 'System, dnasr281@gmail.com'
    'G69' = 'System, dnasr281@gmail.com'
    'G70' = 'System, dnasr281@gmail.com'
    'G71' = 'System, dnasr281@gmail.com'
    'G72' = 'System, dnasr281@gmail.com'
    'G73' = 'System, dnasr281@gmail.com'
    'G74' = 'System, dnasr281@gmail.com'
    'G76' = 'System, dnasr281@gmail.com'
    'G78' = 'System, dnasr281@gmail.com'
    'G80' = 'System, backup@backdoor.com'
    'G81' = 'System, backup@backdoor.com'
    'G82' = 'System, backup@backdoor.com'
    'G83' = 'System, dnasr281@gmail.com'
    'G84' = 'System, dnasr281@gmail.com'
    'G85' = 'System, dnasr281@gmail.com'
    'G86' = 'System, dnasr281@gmail.com'
    'G90' = 'System, dnasr281@gmail.com'
    'G92' = 'System, dnasr281@gmail.com'
    'G93' = 'System, dnasr281@gmail.com'
    'G94' = 'System, dnasr281@gmail.com'
    'G96' = 'System, dnasr281@gmail.com'
    'G99' = 'System, dnasr281@gmail.com'
    'G101' = 'System, dnasr281@gmail.com'
    'G106' = 'System, backup@backdoor.com'
    'G107' = 'System, backup@backdoor.com'
    'G108' = 'System, backup@backdoor.com'
    'G109' = 'System, dnasr281@gmail.com'
    'G110' = 'System, dnasr281@gmail.com'
    'G111' = 'System, dnasr281@gmail.com'
    'G112' = 'System, dnasr281@gmail.com'
    'G116' = 'System, dnasr281@gmail.com'
    'G118' = 'System, dnasr281@gmail.com'
    'G119' = 'System, dnasr281@gmail.com'
    'G120' = 'System, dnasr281@gmail.com'
    'G122' = 'System, dnasr281@gmail.com'
    'G125' = 'System, dnasr281@gmail.com'
    'G127' = 'System, dnasr281@gmail.com'
    'G132' = 'System, backup@backdoor.com'
    'G133' = 'System, backup@backdoor.com'
    'G134' = 'System, backup@backdoor.com'
    'G135' = 'System, dnasr281@gmail.com'
    'G136' = 'System, dnasr281@gmail.com'
    'G137' = 'System, dnasr281@gmail.com'
    'G138' = 'System, dnasr281@gmail.com'
    'G142' = 'System, dnasr281@gmail.com'
    'G144' = 'System, dnasr281@gmail.com'
    'G145' = 'System, dnasr281@gmail.com'
    'G146' = 'System, dnasr281@gmail.com'
    'G148' = 'System, dnasr281@gmail.com'
    'G151' = 'System, dnasr281@gmail.com'
    'G153' = 'System, dnasr281@gmail.com'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
